# The heat-sink purchase line (row 11: "Heatsinks" / "heat transfer" / "Yohan" / $9.62)
# is being removed from Juan's purchases sheet, leaving only the (now blank) priced cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Juan's purchases")

# Fully clear A11 (removes both its value "Heatsinks" and its bold-ish style)
$ws.Range("A11").Clear()

# Clear the remaining contents of that row (purpose/buyer text, price)
$ws.Range("C11:E11").ClearContents()

# Move / record the active selection on the total cell, as saved in the file
$ws.Range("E12").Select()
